# feat: allow creating notes and vocabulary entries with tags
#
# - fix a typo in an existing ENGLISH vocabulary entry
# - append 10 new vocabulary rows to the ENGLISH sheet
# - append 1 new tagged note to the NOTES sheet
# - leave the workbook with ENGLISH as the active/selected sheet

$wb = $excel.ActiveWorkbook

$wsEnglish = $wb.Worksheets.Item("ENGLISH")
$wsNotes   = $wb.Worksheets.Item("NOTES")

# --- fix existing typo: "disentabgle" -> "disentangle" (ENGLISH!A67) ---
$wsEnglish.Cells.Item(67, 1).Value = "disentangle"

# --- new vocabulary entries appended to ENGLISH (rows 115-124) ---
# columns: A=Word B=Definition C=Synonyms D=Antonyms E=Correct answer count F=Created at G=Tags

$newWords = @(
    @{ Row = 115; Word = "resentment";    Definition = $null;                                                              Synonyms = "bitterness";                 Created = "2021-11-18 13:38:55.766952" },
    @{ Row = 116; Word = "coalesce";      Definition = "com together to form one mass or whole";                           Synonyms = "unite";                       Created = "2021-11-18 13:39:52.321676" },
    @{ Row = 117; Word = "timid";         Definition = "easily frightened";                                                Synonyms = $null;                         Created = "2021-11-18 13:40:17.759858" },
    @{ Row = 118; Word = "reverberation"; Definition = "prolongation of a sound; a continuing effect";                     Synonyms = "resonance";                   Created = "2021-11-18 13:41:26.98325"  },
    @{ Row = 119; Word = "imprisonment";  Definition = $null;                                                              Synonyms = "incarceration";               Created = "2021-11-18 13:42:08.650617" },
    @{ Row = 120; Word = "convention";    Definition = $null;                                                              Synonyms = "agreement;custom";            Created = "2021-11-18 13:44:21.614304" },
    @{ Row = 121; Word = "crestfallen";   Definition = $null;                                                              Synonyms = "disappointed;downhearted";    Created = "2021-11-18 13:45:08.281708" },
    @{ Row = 122; Word = "innate";        Definition = $null;                                                              Synonyms = "natural;inborn";              Created = "2021-11-18 13:46:56.18276"  },
    @{ Row = 123; Word = "muddle";        Definition = $null;                                                              Synonyms = "confuse;bewilder";            Created = "2021-11-18 13:47:26.716003" },
    @{ Row = 124; Word = "resolutely";    Definition = "in an admirably purposeful, determined, and unwavering manner";    Synonyms = $null;                         Created = "2021-11-18 13:48:20.311353" }
)

foreach ($entry in $newWords) {
    $r = $entry.Row
    $wsEnglish.Cells.Item($r, 1).Value = $entry.Word
    if ($entry.Definition) {
        $wsEnglish.Cells.Item($r, 2).Value = $entry.Definition
    }
    if ($entry.Synonyms) {
        $wsEnglish.Cells.Item($r, 3).Value = $entry.Synonyms
    }
    $wsEnglish.Cells.Item($r, 5).Value = 0
    $wsEnglish.Cells.Item($r, 6).Value = $entry.Created
}

# --- new tagged note appended to NOTES (row 27) ---
$wsNotes.Cells.Item(27, 1).Value = "The main thing is to keep the main thing the main thing"
$wsNotes.Cells.Item(27, 2).Value = "essentialism"

# --- restore NOTES selection to A1 before switching away from it ---
$wsNotes.Range("A1").Select()

# --- make ENGLISH the active sheet again with M24 selected ---
$wsEnglish.Activate()
$wsEnglish.Range("M24").Select()
